$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 text (Good Morning -> GIT UPDATE)
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active selection on the sheet (E8)
$ws.Range("E8").Select()
